# Update the "data" column (column A) values on Sheet1 so that the
# short numeric-looking date codes (e.g. "0505") are rewritten in the
# "dd/mm" style (e.g. "05/05"), per commit "debits added and funcions optmized".
#
# Mapping of old -> new values:
#   0505 -> 05/05
#   0605 -> 05/06
#   0705 -> 05/07
#   0805 -> 05/08
#   1105 -> 05/11
#   1305 -> 05/13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map each row in column A to its new text value.
$updates = @{
    2  = "05/05"
    3  = "05/05"
    4  = "05/05"
    5  = "05/06"
    6  = "05/06"
    7  = "05/07"
    8  = "05/07"
    9  = "05/07"
    10 = "05/08"
    11 = "05/08"
    12 = "05/08"
    13 = "05/11"
    14 = "05/11"
    15 = "05/13"
    16 = "05/13"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $updates[$row]
}
